# "Full modification of the flow"
#
# The rule-concept-map sheet goes from 4 rows of "model" rows (CRE PD Model /
# SME PD Model / Auto LGD Model / HE EAD Model, each tagged with the same
# long numeric id list) down to 2 rows describing actual rules (Rule 1 / Rule
# 2) tied to CCAR ids and short comma lists.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the old row 3 (Auto LGD Model) and row 4 (Auto LDG Model) entirely -
# the sheet shrinks from A1:D4 to A1:D2.
$ws.Rows("3:4").Delete()

# Row 1: 34/CRE PD Model/CRE PD Model/<id list>  ->  76/CCAR1/Rule 1 is rule 1/1,2
$ws.Range("A1").Value = 76
$ws.Range("B1").Value = "CCAR1"
$ws.Range("C1").Value = "Rule 1 is rule 1"
$ws.Range("D1").Value = "1,2"

# Row 2: 38/SME PD Model/SME PD Model/<id list>  ->  2/CCAR2/Rule 2 is rule 2/3,4
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "CCAR2"
$ws.Range("C2").Value = "Rule 2 is rule 2"
$ws.Range("D2").Value = "3,4"

# Column D no longer auto-fits the long id-list text; it gets a plain fixed
# width instead (was width="41.5703125" bestFit="1").
$ws.Columns("D").ColumnWidth = 22

# The stale "I13" selection (left over from editing far outside the old
# table) is cleared back to the top-left cell along with everything else.
$ws.Range("A1").Select() | Out-Null
